$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 747.5
$ws.Range("I2").Value = 262.33334
$ws.Range("K2").Value = 262.33334
$ws.Range("M2").Value = -149.33334
$ws.Range("H19").Value = 2559.1875
$ws.Range("I19").Value = 2076.8462
$ws.Range("J19").Value = 2889.2104
$ws.Range("K19").Value = 2076.8462
$ws.Range("L19").Value = 2889.2104
$ws.Range("M19").Value = -1901.8462
$ws.Range("N19").Value = -3239.2104
$ws.Range("H29").Value = 7886.8887
$ws.Range("J29").Value = 9497.143
$ws.Range("L29").Value = 28491.429
$ws.Range("N29").Value = -29053.429
$ws.Range("H86").Value = 3558.9375
$ws.Range("I86").Value = 3196.2
$ws.Range("J86").Value = 9000
$ws.Range("K86").Value = 3196.2
$ws.Range("L86").Value = 9000
$ws.Range("M86").Value = -2073.2
$ws.Range("N86").Value = -11246
$ws.Range("H89").Value = 3558.9375
$ws.Range("I89").Value = 3196.2
$ws.Range("J89").Value = 9000
$ws.Range("K89").Value = 15981
$ws.Range("L89").Value = 45000
$ws.Range("M89").Value = -10365
$ws.Range("N89").Value = -56232
$ws.Range("H92").Value = 984.94446
$ws.Range("J92").Value = 1043
$ws.Range("L92").Value = 1043
$ws.Range("N92").Value = -3539
$ws.Range("H112").Value = 5146.119
$ws.Range("J112").Value = 5169.2197
$ws.Range("L112").Value = 15507.6591
$ws.Range("N112").Value = -17723.6591
$ws.Range("H116").Value = 4206.0713
$ws.Range("I116").Value = 3520.2222
$ws.Range("J116").Value = 5440.6
$ws.Range("K116").Value = 3520.2222
$ws.Range("L116").Value = 5440.6
$ws.Range("M116").Value = -78.22220000000016
$ws.Range("N116").Value = -12324.6
$ws.Range("H138").Value = 3153.034
$ws.Range("J138").Value = 3981.1316
$ws.Range("L138").Value = 11943.3948
$ws.Range("N138").Value = -22223.3948

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1645.7333
$ws.Range("I88").Value = 1424.75
$ws.Range("J88").Value = 1726.091
$ws.Range("K88").Value = 1424.75
$ws.Range("L88").Value = 1726.091
$ws.Range("M88").Value = -1018.75
$ws.Range("N88").Value = -2538.091
$ws.Range("H91").Value = 1645.7333
$ws.Range("I91").Value = 1424.75
$ws.Range("J91").Value = 1726.091
$ws.Range("K91").Value = 1424.75
$ws.Range("L91").Value = 1726.091
$ws.Range("M91").Value = -20.75
$ws.Range("N91").Value = -4534.091
$ws.Range("H98").Value = 73845.92
$ws.Range("J98").Value = 73845.92
$ws.Range("L98").Value = 73845.92
$ws.Range("N98").Value = -79835.92
$ws.Range("H102").Value = 1836.5714
$ws.Range("I102").Value = 1785.5385
$ws.Range("K102").Value = 1785.5385
$ws.Range("M102").Value = -163.5385000000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2405.138
$ws.Range("I107").Value = 2392.08
$ws.Range("K107").Value = 2392.08
$ws.Range("M107").Value = -472.0799999999999
$ws.Range("H134").Value = 4360.5
$ws.Range("I134").Value = 3163.818
$ws.Range("K134").Value = 9491.454000000002
$ws.Range("M134").Value = -6956.454000000002

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 80233.28999999999
$ws.Range("I92").Value = 15000
$ws.Range("K92").Value = 15000
$ws.Range("M92").Value = -12504
$ws.Range("H107").Value = 21740444
$ws.Range("I107").Value = 33334084
$ws.Range("J107").Value = 2372
$ws.Range("K107").Value = 33334084
$ws.Range("L107").Value = 2372
$ws.Range("M107").Value = -33332164
$ws.Range("N107").Value = -6212
$ws.Range("H122").Value = 2223.125
$ws.Range("I122").Value = 326.42856
$ws.Range("K122").Value = 979.28568
$ws.Range("M122").Value = 1470.71432
$ws.Range("H134").Value = 1039.9678
$ws.Range("I134").Value = 897.8889
$ws.Range("K134").Value = 2693.6667
$ws.Range("M134").Value = -158.6667000000002

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 860.6
$ws.Range("I69").Value = 860.6
$ws.Range("K69").Value = 2581.8
$ws.Range("M69").Value = -1770.8
$ws.Range("H70").Value = 5477.125
$ws.Range("I70").Value = 5477.125
$ws.Range("K70").Value = 16431.375
$ws.Range("M70").Value = -16116.375
$ws.Range("H72").Value = 860.6
$ws.Range("I72").Value = 860.6
$ws.Range("K72").Value = 7745.400000000001
$ws.Range("M72").Value = -3689.400000000001
$ws.Range("H73").Value = 5477.125
$ws.Range("I73").Value = 5477.125
$ws.Range("K73").Value = 16431.375
$ws.Range("M73").Value = -15339.375
$ws.Range("H94").Value = 17428.285
$ws.Range("J94").Value = 17428.285
$ws.Range("L94").Value = 52284.855
$ws.Range("N94").Value = -53636.855
$ws.Range("H131").Value = 2772.0476
$ws.Range("I131").Value = 746.375
$ws.Range("J131").Value = 4018.6155
$ws.Range("K131").Value = 2239.125
$ws.Range("L131").Value = 12055.8465
$ws.Range("M131").Value = 2800.875
$ws.Range("N131").Value = -22135.8465
$ws.Range("H133").Value = 13560.889
$ws.Range("I133").Value = 10879.167
$ws.Range("J133").Value = 14901.75
$ws.Range("K133").Value = 32637.501
$ws.Range("L133").Value = 44705.25
$ws.Range("M133").Value = -27577.501
$ws.Range("N133").Value = -54825.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 21654.5
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("H80").Value = 4832.5625
$ws.Range("I80").Value = 4746.5
$ws.Range("J80").Value = 4918.625
$ws.Range("K80").Value = 4746.5
$ws.Range("L80").Value = 4918.625
$ws.Range("M80").Value = -3748.5
$ws.Range("N80").Value = -6914.625
$ws.Range("H83").Value = 4832.5625
$ws.Range("I83").Value = 4746.5
$ws.Range("J83").Value = 4918.625
$ws.Range("K83").Value = 23732.5
$ws.Range("L83").Value = 24593.125
$ws.Range("M83").Value = -18740.5
$ws.Range("N83").Value = -34577.125
$ws.Range("H133").Value = 88431.414
$ws.Range("J133").Value = 88431.414
$ws.Range("L133").Value = 88431.414
$ws.Range("N133").Value = -98551.414
$ws.Range("N58").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1850.6666
$ws.Range("I82").Value = 1459.6666
$ws.Range("J82").Value = 2632.6667
$ws.Range("K82").Value = 1459.6666
$ws.Range("L82").Value = 2632.6667
$ws.Range("M82").Value = -1098.6666
$ws.Range("N82").Value = -3354.6667
$ws.Range("H85").Value = 1850.6666
$ws.Range("I85").Value = 1459.6666
$ws.Range("J85").Value = 2632.6667
$ws.Range("K85").Value = 1459.6666
$ws.Range("L85").Value = 2632.6667
$ws.Range("M85").Value = -211.6666
$ws.Range("N85").Value = -5128.6667
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("H99").Value = 71828
$ws.Range("I99").Value = 17199
$ws.Range("J99").Value = 99142.5
$ws.Range("K99").Value = 17199
$ws.Range("L99").Value = 99142.5
$ws.Range("M99").Value = -14204
$ws.Range("N99").Value = -105132.5
$ws.Range("H132").Value = 6754.8
$ws.Range("I132").Value = 6915.9287
$ws.Range("J132").Value = 6549.727
$ws.Range("K132").Value = 20747.7861
$ws.Range("L132").Value = 19649.181
$ws.Range("M132").Value = -18217.7861
$ws.Range("N132").Value = -24709.181
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6764.125
$ws.Range("I62").Value = 5827.25
$ws.Range("J62").Value = 7701
$ws.Range("K62").Value = 5827.25
$ws.Range("L62").Value = 7701
$ws.Range("M62").Value = -5203.25
$ws.Range("N62").Value = -8949
$ws.Range("H65").Value = 6764.125
$ws.Range("I65").Value = 5827.25
$ws.Range("J65").Value = 7701
$ws.Range("K65").Value = 29136.25
$ws.Range("L65").Value = 38505
$ws.Range("M65").Value = -26016.25
$ws.Range("N65").Value = -44745
$ws.Range("H74").Value = 5773.143
$ws.Range("J74").Value = 5755.3335
$ws.Range("L74").Value = 5755.3335
$ws.Range("N74").Value = -7627.3335
$ws.Range("H77").Value = 5773.143
$ws.Range("J77").Value = 5755.3335
$ws.Range("L77").Value = 17266.0005
$ws.Range("N77").Value = -26626.0005
